# ALT_02 - Modificacao no Diagrama de Use Case e nas suas Especificacoes (Parcial)
#
# Changes applied (per unified diff):
#  1. Move the hidden "_GoBack" bookmark from the paragraph that ends with
#     'clicar em "Cadastrar".' (inside GATILHO) up to the empty, centered
#     paragraph right before the table (the 2nd paragraph in the document).
#  2. Split the run containing "RF001" into two runs: "RF00" and "2".
#  3. Delete the last table row ("REGRAS DE NEGOCIO: RE001, RE002, RE003").

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Step 1a: remove the bookmark from the GATILHO / "Cadastrar" paragraph ---
$gatilhoIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*GATILHO*" -and $t -like "*Cadastrar*") {
        $gatilhoIdx = $i
        break
    }
}

if ($gatilhoIdx -ne -1) {
    $p = $d.Paragraphs($gatilhoIdx)
    $r = $p.Range
    $body = '<w:p w:rsidR="006F0C48" w:rsidRDefault="006F0C48" w:rsidP="007F5348"><w:r w:rsidRPr="006F0C48"><w:rPr><w:b/></w:rPr><w:t>GATILHO:</w:t></w:r><w:r w:rsidR="00A57170"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00A57170"><w:t xml:space="preserve">Ao abrir o </w:t></w:r><w:r w:rsidR="00F155FC"><w:t>site</w:t></w:r><w:r w:rsidR="00A57170"><w:t xml:space="preserve"> o ator ir&#225; </w:t></w:r><w:r w:rsidR="00DA053F"><w:t>clicar em &#8220;Cadastrar&#8221;.</w:t></w:r></w:p>'
    [void]$r.InsertXML($pkgHeader + $body + $pkgFooter)
}

# --- Step 1b: add the bookmark to the empty centered paragraph (2nd paragraph) ---
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$body2 = '<w:p w:rsidR="005A0BCD" w:rsidRPr="005A0BCD" w:rsidRDefault="005A0BCD" w:rsidP="005A0BCD"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
[void]$r2.InsertXML($pkgHeader + $body2 + $pkgFooter)

# --- Step 2: split "RF001" into "RF00" + "2" (two runs) ---
$rfIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*RF001*") {
        $rfIdx = $i
        break
    }
}

if ($rfIdx -ne -1) {
    $p3 = $d.Paragraphs($rfIdx)
    $r3 = $p3.Range
    $body3 = '<w:p w:rsidR="006F0C48" w:rsidRPr="00A57170" w:rsidRDefault="00A641F6" w:rsidP="006F0C48"><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">ID: </w:t></w:r><w:r><w:t>RF00</w:t></w:r><w:r><w:t>2</w:t></w:r></w:p>'
    [void]$r3.InsertXML($pkgHeader + $body3 + $pkgFooter)
}

# --- Step 3: delete the last table row (REGRAS DE NEGOCIO: RE001, RE002, RE003) ---
$tbl = $d.Tables(1)
$lastRow = $tbl.Rows($tbl.Rows.Count)
if ($lastRow.Range.Text -like "*REGRAS DE NEG*") {
    [void]$lastRow.Delete()
}

Write-Host "edit complete"
